$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 151 - shifts existing rows 151:199 down to 152:200
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new data record
$ws.Range("A151").Value = 10
$ws.Range("B151").Value = "Vega Modelo de Temuco"
$ws.Range("C151").Value = "La Araucanía"
$ws.Range("D151").Value = 44559
$ws.Range("E151").Value = 9
$ws.Range("F151").Value = "Fruta"
$ws.Range("G151").Value = 100103
$ws.Range("H151").Value = "Frutos de hueso (carozo)"
$ws.Range("I151").Value = 100103001
$ws.Range("J151").Value = "Cereza"
$ws.Range("K151").Value = "Lapins"
$ws.Range("L151").Value = "Primera"
$ws.Range("M151").Value = 1170
$ws.Range("N151").Value = 5000
$ws.Range("O151").Value = 5500
$ws.Range("P151").Value = 5222
$ws.Range("Q151").Value = "$/bandeja 10 kilos"
$ws.Range("R151").Value = "Región del Maule"
$ws.Range("S151").Value = 522
$ws.Range("T151").Value = 10

# Apply the same date style (numFmt) used by the rest of column D to the new row's date cell
$ws.Range("D151").NumberFormat = $ws.Range("D152").NumberFormat
